# Update "想去人数" (number of people interested) values on the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets to match the
# newly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2150
$ws1.Range("F4").Value  = 38
$ws1.Range("F5").Value  = 11312
$ws1.Range("F7").Value  = 315
$ws1.Range("F9").Value  = 11250
$ws1.Range("F10").Value = 457
$ws1.Range("F12").Value = 62
$ws1.Range("F13").Value = 1736
$ws1.Range("F14").Value = 5617
$ws1.Range("F15").Value = 101
$ws1.Range("F16").Value = 3464

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 2150
$ws4.Range("F5").Value  = 38
$ws4.Range("F7").Value  = 11312
$ws4.Range("F9").Value  = 315
$ws4.Range("F11").Value = 11250
$ws4.Range("F12").Value = 457
$ws4.Range("F14").Value = 62
$ws4.Range("F15").Value = 1736
$ws4.Range("F16").Value = 5617
$ws4.Range("F17").Value = 101
$ws4.Range("F18").Value = 3464

$wb.Save()
